# Auto-generated script to apply scheduled market-data / profit refresh
# to the Kraken_Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# For each touched row, the live market price / computed profit columns
# (H:currentAveragePrice, I:currentAveragePriceNQ, J:currentAveragePriceHQ,
#  K:LevePriceNQ, L:LevePriceHQ, M:LeveProfitNQ, N:LeveProfitHQ) are refreshed.
# A $null assignment clears/removes a cell whose profit became undefined.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
# Row 2
$ws.Range("H2").Value = 112.166664
$ws.Range("I2").Value = 104.6
$ws.Range("J2").Value = 150
$ws.Range("K2").Value = 104.6
$ws.Range("L2").Value = 150
$ws.Range("M2").Value = 8.400000000000006
$ws.Range("N2").Value = -376
# Row 18
$ws.Range("H18").Value = 22042.857
$ws.Range("I18").Value = 25980
$ws.Range("K18").Value = 25980
$ws.Range("M18").Value = -25696
# Row 40
$ws.Range("H40").Value = 7199.6
$ws.Range("J40").Value = 7199.6
$ws.Range("L40").Value = 7199.6
$ws.Range("N40").Value = -7549.6
# Row 76
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").Value = $null
# Row 79
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").Value = $null
# Row 86
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = $null
$ws.Range("N86").Value = 0
# Row 89
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = $null
$ws.Range("N89").Value = 0
# Row 106
$ws.Range("H106").Value = 2856.4285
$ws.Range("I106").Value = 1995
$ws.Range("K106").Value = 1995
$ws.Range("M106").Value = -1364
# Row 125
$ws.Range("H125").Value = 670.4
$ws.Range("I125").Value = 638
$ws.Range("K125").Value = 5742
$ws.Range("M125").Value = -3282
# Row 137
$ws.Range("H137").Value = 2337.125
$ws.Range("I137").Value = 2314.1428
$ws.Range("J137").Value = 2498
$ws.Range("K137").Value = 6942.428400000001
$ws.Range("L137").Value = 7494
$ws.Range("M137").Value = -4392.428400000001
$ws.Range("N137").Value = -12594

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
# Row 61
$ws.Range("H61").Value = 2499
$ws.Range("I61").Value = 2499.5
$ws.Range("J61").Value = 2498
$ws.Range("K61").Value = 2499.5
$ws.Range("L61").Value = 2498
$ws.Range("M61").Value = -2287.5
$ws.Range("N61").Value = -2922
# Row 63
$ws.Range("H63").Value = 5000
$ws.Range("I63").Value = 5000
$ws.Range("K63").Value = 5000
$ws.Range("M63").Value = -4314
# Row 66
$ws.Range("H66").Value = 5000
$ws.Range("I66").Value = 5000
$ws.Range("K66").Value = 25000
$ws.Range("M66").Value = -21568
# Row 88
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = $null
$ws.Range("M88").Value = $null
$ws.Range("N88").Value = 0
# Row 91
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = $null
$ws.Range("M91").Value = $null
$ws.Range("N91").Value = 0
# Row 132
$ws.Range("H132").Value = 1168.6666
$ws.Range("I132").Value = 1156.7273
$ws.Range("K132").Value = 3470.1819
$ws.Range("M132").Value = -940.1819
# Row 136
$ws.Range("H136").Value = 2499
$ws.Range("I136").Value = 2499.5
$ws.Range("J136").Value = 2498
$ws.Range("K136").Value = 7498.5
$ws.Range("L136").Value = 7494
$ws.Range("M136").Value = -4948.5
$ws.Range("N136").Value = -12594

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
# Row 105
$ws.Range("H105").Value = 2699.7273
$ws.Range("J105").Value = 2633.3333
$ws.Range("L105").Value = 2633.3333
$ws.Range("N105").Value = -6127.3333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
# Row 31
$ws.Range("H31").Value = 2722.5417
$ws.Range("I31").Value = 1384.4
$ws.Range("K31").Value = 1384.4
$ws.Range("M31").Value = -1089.4
# Row 34
$ws.Range("H34").Value = 2722.5417
$ws.Range("I34").Value = 1384.4
$ws.Range("K34").Value = 1384.4
$ws.Range("M34").Value = -1182.4
# Row 35
$ws.Range("H35").Value = 824.5
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").Value = $null
# Row 58
$ws.Range("H58").Value = 500
$ws.Range("I58").Value = 500
$ws.Range("K58").Value = 500
$ws.Range("M58").Value = -297
# Row 62
$ws.Range("H62").Value = 20000
$ws.Range("I62").Value = 20000
$ws.Range("K62").Value = 20000
$ws.Range("M62").Value = -19376
# Row 65
$ws.Range("H65").Value = 20000
$ws.Range("I65").Value = 20000
$ws.Range("K65").Value = 100000
$ws.Range("M65").Value = -96880
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = $null
$ws.Range("N87").Value = 0
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = $null
$ws.Range("N90").Value = 0
# Row 132
$ws.Range("H132").Value = 2512.2222
$ws.Range("I132").Value = 2201.25
$ws.Range("K132").Value = 6603.75
$ws.Range("M132").Value = -4073.75
# Row 136
$ws.Range("H136").Value = 500
$ws.Range("I136").Value = 500
$ws.Range("K136").Value = 1500
$ws.Range("M136").Value = 1050

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)
# Row 87
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").Value = $null
# Row 90
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").Value = $null
# Row 107
$ws.Range("H107").Value = 300
$ws.Range("I107").Value = 300
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 900
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = $null
$ws.Range("N107").Value = 1020

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
# Row 80
$ws.Range("H80").Value = 68664.336
$ws.Range("I80").Value = 2993.5
$ws.Range("J80").Value = 200006
$ws.Range("K80").Value = 2993.5
$ws.Range("L80").Value = 200006
$ws.Range("M80").Value = -1995.5
$ws.Range("N80").Value = -202002
# Row 83
$ws.Range("H83").Value = 68664.336
$ws.Range("I83").Value = 2993.5
$ws.Range("J83").Value = 200006
$ws.Range("K83").Value = 14967.5
$ws.Range("L83").Value = 1000030
$ws.Range("M83").Value = -9975.5
$ws.Range("N83").Value = -1010014
# Row 132
$ws.Range("H132").Value = 2261.2666
$ws.Range("I132").Value = 1784.9166
$ws.Range("K132").Value = 5354.7498
$ws.Range("M132").Value = -2824.7498
# Row 135
$ws.Range("H135").Value = 90000
$ws.Range("J135").Value = 90000
$ws.Range("L135").Value = 90000
$ws.Range("N135").Value = -100140

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
# Row 46
$ws.Range("H46").Value = 2935.1765
$ws.Range("J46").Value = 3249.8572
$ws.Range("L46").Value = 3249.8572
$ws.Range("N46").Value = -3625.8572
# Row 68
$ws.Range("H68").Value = 2566.5557
$ws.Range("I68").Value = 2512.375
$ws.Range("K68").Value = 2512.375
$ws.Range("M68").Value = -1763.375
# Row 71
$ws.Range("H71").Value = 2566.5557
$ws.Range("I71").Value = 2512.375
$ws.Range("K71").Value = 12561.875
$ws.Range("M71").Value = -8817.875
# Row 82
$ws.Range("H82").Value = 1991.125
$ws.Range("I82").Value = 1990
$ws.Range("J82").Value = 1993
$ws.Range("K82").Value = 1990
$ws.Range("L82").Value = 1993
$ws.Range("M82").Value = -1629
$ws.Range("N82").Value = -2715
# Row 85
$ws.Range("H85").Value = 1991.125
$ws.Range("I85").Value = 1990
$ws.Range("J85").Value = 1993
$ws.Range("K85").Value = 1990
$ws.Range("L85").Value = 1993
$ws.Range("M85").Value = -742
$ws.Range("N85").Value = -4489
# Row 132
$ws.Range("H132").Value = 7333
$ws.Range("I132").Value = 9499.5
$ws.Range("K132").Value = 28498.5
$ws.Range("M132").Value = -25968.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
# Row 46
$ws.Range("H46").Value = 99995
$ws.Range("J46").Value = 99995
$ws.Range("L46").Value = 99995
$ws.Range("N46").Value = -100457
# Row 69
$ws.Range("H69").Value = 7180.3335
$ws.Range("J69").Value = 7180.3335
$ws.Range("L69").Value = 7180.3335
$ws.Range("N69").Value = -8678.333500000001
# Row 72
$ws.Range("H72").Value = 7180.3335
$ws.Range("J72").Value = 7180.3335
$ws.Range("L72").Value = 21541.0005
$ws.Range("N72").Value = -29029.0005
# Row 132
$ws.Range("H132").Value = 3165
$ws.Range("I132").Value = 3098
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 9294
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -6764
$ws.Range("N132").Value = -15560
# Row 134
$ws.Range("H134").Value = 99995
$ws.Range("J134").Value = 99995
$ws.Range("L134").Value = 299985
$ws.Range("N134").Value = -305055
